$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false,
                             $true, 1, $false, $replace, 2)
}

Replace-Text "852×3=" "729×9="
Replace-Text "298×2=" "583×2="
Replace-Text "154×5=" "976×8="
Replace-Text "796×3=" "722×8="
Replace-Text "606×6=" "944×8="
Replace-Text "475×8=" "864×9="
Replace-Text "945×8=" "363×3="
Replace-Text "844×8=" "437×6="
Replace-Text "289×7=" "256×5="
Replace-Text "332×9=" "765×5="
Replace-Text "320×5=" "864×4="
Replace-Text "830×4=" "153×9="
Replace-Text "354×9=" "146×4="
Replace-Text "563×8=" "372×3="
Replace-Text "788×6=" "331×3="
Replace-Text "130×4=" "117×6="
Replace-Text "271×2=" "757×7="
Replace-Text "728×8=" "317×7="
Replace-Text "633×3=" "964×6="
Replace-Text "891×9=" "631×9="
Replace-Text "430×8=" "915×7="
Replace-Text "647×6=" "468×9="
Replace-Text "756×3=" "352×3="
Replace-Text "400×7=" "247×4="
Replace-Text "648×4=" "258×7="
